$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values in A2/A3, Q2/Q3, R2/R3 using Value2 to avoid Variant wrapper issues
$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$ws.Range("A2").Value2 = $a3
$ws.Range("A3").Value2 = $a2

$q2 = $ws.Range("Q2").Value2
$q3 = $ws.Range("Q3").Value2
$ws.Range("Q2").Value2 = $q3
$ws.Range("Q3").Value2 = $q2

$r2 = $ws.Range("R2").Value2
$r3 = $ws.Range("R3").Value2
$ws.Range("R2").Value2 = $r3
$ws.Range("R3").Value2 = $r2
